# Update LR-pair TPM results with recomputed values and two new
# target-cluster rows (FAPs, MuSCs) for every sending cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Gp6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.73076433333334
$ws.Range("H2").Value = 53.19229300000001
$ws.Range("I2").Value = 0.004631884691211661
$ws.Range("J2").Value = 0.00463188469121166
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.053317
$ws.Range("N2").Value = 0.159951
$ws.Range("O2").Value = 0.5209400637695697
$ws.Range("P2").Value = 0.5209400637695697
$ws.Range("Q2").Value = 0.9453511619603335
$ws.Range("R2").Value = 8.508160457643001
$ws.Range("S2").Value = 0.002412934306413096
$ws.Range("T2").Value = 0.002412934306413096

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Gp6"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.73076433333334
$ws.Range("H3").Value = 53.19229300000001
$ws.Range("I3").Value = 0.004631884691211661
$ws.Range("J3").Value = 0.00463188469121166
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04903066666666667
$ws.Range("N3").Value = 0.147092
$ws.Range("O3").Value = 0.4790599362304302
$ws.Range("P3").Value = 0.4790599362304302
$ws.Range("Q3").Value = 0.8693511957728891
$ws.Range("R3").Value = 7.824160761956001
$ws.Range("S3").Value = 0.002218950384798564
$ws.Range("T3").Value = 0.002218950384798564

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Gp6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3580.644531333333
$ws.Range("H4").Value = 10741.933594
$ws.Range("I4").Value = 0.9353873458333681
$ws.Range("J4").Value = 0.935387345833368
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.053317
$ws.Range("N4").Value = 0.159951
$ws.Range("O4").Value = 0.5209400637695697
$ws.Range("P4").Value = 0.5209400637695697
$ws.Range("Q4").Value = 190.9092244770993
$ws.Range("R4").Value = 1718.183020293894
$ws.Range("S4").Value = 0.4872807435876834
$ws.Range("T4").Value = 0.4872807435876833

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Gp6"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3580.644531333333
$ws.Range("H5").Value = 10741.933594
$ws.Range("I5").Value = 0.9353873458333681
$ws.Range("J5").Value = 0.935387345833368
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04903066666666667
$ws.Range("N5").Value = 0.147092
$ws.Range("O5").Value = 0.4790599362304302
$ws.Range("P5").Value = 0.4790599362304302
$ws.Range("Q5").Value = 175.5613884676276
$ws.Range("R5").Value = 1580.052496208648
$ws.Range("S5").Value = 0.4481066022456848
$ws.Range("T5").Value = 0.4481066022456847

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Gp6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 227.2177583333333
$ws.Range("H6").Value = 681.653275
$ws.Range("I6").Value = 0.0593570833501536
$ws.Range("J6").Value = 0.05935708335015359
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.053317
$ws.Range("N6").Value = 0.159951
$ws.Range("O6").Value = 0.5209400637695697
$ws.Range("P6").Value = 0.5209400637695697
$ws.Range("Q6").Value = 12.11456922105833
$ws.Range("R6").Value = 109.031122989525
$ws.Range("S6").Value = 0.03092148278560468
$ws.Range("T6").Value = 0.03092148278560468

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Gp6"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 227.2177583333333
$ws.Range("H7").Value = 681.653275
$ws.Range("I7").Value = 0.0593570833501536
$ws.Range("J7").Value = 0.05935708335015359
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.04903066666666667
$ws.Range("N7").Value = 0.147092
$ws.Range("O7").Value = 0.4790599362304302
$ws.Range("P7").Value = 0.4790599362304302
$ws.Range("Q7").Value = 11.14063816958889
$ws.Range("R7").Value = 100.2657435263
$ws.Range("S7").Value = 0.02843560056454892
$ws.Range("T7").Value = 0.02843560056454891

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Col1a2"
$ws.Range("C8").Value = "Gp6"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.387458333333333
$ws.Range("H8").Value = 7.162374999999999
$ws.Range("I8").Value = 0.0006236861252666267
$ws.Range("J8").Value = 0.0006236861252666266
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.053317
$ws.Range("N8").Value = 0.159951
$ws.Range("O8").Value = 0.5209400637695697
$ws.Range("P8").Value = 0.5209400637695697
$ws.Range("Q8").Value = 0.1272921159583333
$ws.Range("R8").Value = 1.145629043625
$ws.Range("S8").Value = 0.0003249030898685924
$ws.Range("T8").Value = 0.0003249030898685923

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Col1a2"
$ws.Range("C9").Value = "Gp6"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.387458333333333
$ws.Range("H9").Value = 7.162374999999999
$ws.Range("I9").Value = 0.0006236861252666267
$ws.Range("J9").Value = 0.0006236861252666266
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.04903066666666667
$ws.Range("N9").Value = 0.147092
$ws.Range("O9").Value = 0.4790599362304302
$ws.Range("P9").Value = 0.4790599362304302
$ws.Range("Q9").Value = 0.1170586737222222
$ws.Range("R9").Value = 1.0535280635
$ws.Range("S9").Value = 0.0002987830353980343
$ws.Range("T9").Value = 0.0002987830353980343

